$wb = $excel.ActiveWorkbook

# Rename the export tab from "BNoEVC" to "BNoGP" to match the csv file / var name
$wsData = $wb.Worksheets.Item("BNoEVC")
$wsData.Name = "BNoGP"

$wsAbout = $wb.Worksheets.Item("About")

# Update the "About" sheet title cell (A1) to reference the new name/description
$wsAbout.Range("A1").Value = "BNoEVC BAU Number of Gas Pumps"

# Update selections to match the saved view state
$wsAbout.Range("B7").Select()

# Make the data sheet the active tab, with its own saved selection
$wsData.Activate()
$wsData.Range("E8").Select()
